$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.406.64"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.574.55"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.73"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.82"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.247"
$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.799.48"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.565.34"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.424.75"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.72"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.29"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0687"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.98"
$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("E25").Value = "  -1.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.11"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.95"

$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("E31").Value = "  +3.50%  "

$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("E33").Value = "  -0.54%  "

$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.382.44"
$ws.Range("E35").Value = "  -1.14%  "

$ws.Range("E36").Value = "  +4.32%  "

$ws.Range("E37").Value = "  -2.80%  "

$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("E40").Value = "  -1.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.522"
$ws.Range("E41").Value = "  -2.21%  "

$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("E46").Value = "  -4.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.48"
$ws.Range("E47").Value = "  -0.78%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.919"
$ws.Range("E48").Value = "  -6.11%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.711.31"
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "85.46"
$ws.Range("E51").Value = "  -0.74%  "
